$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040980325599067
$ws.Range("D2").Value = 1.040951545586742
$ws.Range("E2").Value = 1.044637614288245
$ws.Range("F2").Value = 1.039742532480231
$ws.Range("I2").Value = 1.033366393062693
$ws.Range("J2").Value = 1.046064003737711
$ws.Range("K2").Value = 1.043732438570256
$ws.Range("L2").Value = 1.047408103293208
$ws.Range("M2").Value = 1.042526855358008
$ws.Range("N2").Value = 1.047549534608237
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.043560390678937
$ws.Range("D3").Value = 1.042931612049246
$ws.Range("E3").Value = 1.047169698320997
$ws.Range("F3").Value = 1.042971897928034
$ws.Range("I3").Value = 1.034014529080863
$ws.Range("J3").Value = 1.048281792677045
$ws.Range("K3").Value = 1.045519962940586
$ws.Range("L3").Value = 1.049746984909431
$ws.Range("M3").Value = 1.045560143188544
$ws.Range("N3").Value = 1.049770473062249
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045217827346884
$ws.Range("D4").Value = 1.044202284581959
$ws.Range("E4").Value = 1.04879654112603
$ws.Range("F4").Value = 1.045048007310766
$ws.Range("I4").Value = 1.034427767191165
$ws.Range("J4").Value = 1.049704967624586
$ws.Range("K4").Value = 1.046665681990339
$ws.Range("L4").Value = 1.05124856967063
$ws.Range("M4").Value = 1.047509303829328
$ws.Range("N4").Value = 1.051195669081456
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045911804934407
$ws.Range("D5").Value = 1.044734002171011
$ws.Range("E5").Value = 1.049477762949249
$ws.Range("F5").Value = 1.045917666173718
$ws.Range("I5").Value = 1.034600041902576
$ws.Range("J5").Value = 1.050300490548349
$ws.Range("K5").Value = 1.047144779149495
$ws.Range("L5").Value = 1.051877071475913
$ws.Range("M5").Value = 1.048325571073269
$ws.Range("N5").Value = 1.051792037716078
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046028164216777
$ws.Range("D6").Value = 1.044823136465233
$ws.Range("E6").Value = 1.049591986618544
$ws.Range("F6").Value = 1.046063504490178
$ws.Range("I6").Value = 1.034628883140018
$ws.Range("J6").Value = 1.050400320322564
$ws.Range("K6").Value = 1.047225072896249
$ws.Range("L6").Value = 1.051982439518074
$ws.Range("M6").Value = 1.048462443252607
$ws.Range("N6").Value = 1.051892009260023
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045227111249204
$ws.Range("D7").Value = 1.044209399062477
$ws.Range("E7").Value = 1.048805654172493
$ws.Range("F7").Value = 1.045059639936707
$ws.Range("I7").Value = 1.034430074805532
$ws.Range("J7").Value = 1.049712935867878
$ws.Range("K7").Value = 1.046672093710864
$ws.Range("L7").Value = 1.051256978515828
$ws.Range("M7").Value = 1.047520223128049
$ws.Range("N7").Value = 1.051203648640567
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041854807709194
$ws.Range("D8").Value = 1.041622939285645
$ws.Range("E8").Value = 1.045495785598999
$ws.Range("F8").Value = 1.040836766488335
$ws.Range("I8").Value = 1.033586720022479
$ws.Range("J8").Value = 1.046816015826105
$ws.Range("K8").Value = 1.044338836069791
$ws.Range("L8").Value = 1.048201029271765
$ws.Range("M8").Value = 1.043554836301965
$ws.Range("N8").Value = 1.048302614640054
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035816834899227
$ws.Range("D9").Value = 1.036981820917625
$ws.Range("E9").Value = 1.039571419036545
$ws.Range("F9").Value = 1.033287696457957
$ws.Range("I9").Value = 1.032052537967777
$ws.Range("J9").Value = 1.041617291600494
$ws.Range("K9").Value = 1.040141199192216
$ws.Range("L9").Value = 1.042722390232306
$ws.Range("M9").Value = 1.036459147578809
$ws.Range("N9").Value = 1.04309650762974
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031722513114367
$ws.Range("D10").Value = 1.033828007252485
$ws.Range("E10").Value = 1.035555414806415
$ws.Range("F10").Value = 1.028176232477886
$ws.Range("I10").Value = 1.030996012447559
$ws.Range("J10").Value = 1.038084050366391
$ws.Range("K10").Value = 1.037281384092493
$ws.Range("L10").Value = 1.03900260258369
$ws.Range("M10").Value = 1.031650019340162
$ws.Range("N10").Value = 1.039558248787816
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02993217102625
$ws.Range("D11").Value = 1.032447374015807
$ws.Range("E11").Value = 1.033799635309964
$ws.Range("F11").Value = 1.025942806091055
$ws.Range("I11").Value = 1.03053020102539
$ws.Range("J11").Value = 1.036537158290368
$ws.Range("K11").Value = 1.036027697905854
$ws.Range("L11").Value = 1.037374915196527
$ws.Range("M11").Value = 1.029547598193605
$ws.Range("N11").Value = 1.038009159947614
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029264437584253
$ws.Range("D12").Value = 1.031932217360918
$ws.Range("E12").Value = 1.033144841049689
$ws.Range("F12").Value = 1.025110062754896
$ws.Range("I12").Value = 1.030355898267816
$ws.Range("J12").Value = 1.035959938464535
$ws.Range("K12").Value = 1.035559644575437
$ws.Range("L12").Value = 1.036767678524813
$ws.Range("M12").Value = 1.02876353607013
$ws.Range("N12").Value = 1.03743112040342
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029407793353934
$ws.Range("D13").Value = 1.032042826750649
$ws.Range("E13").Value = 1.033285416691175
$ws.Range("F13").Value = 1.025288833593098
$ws.Range("I13").Value = 1.030393345172306
$ws.Range("J13").Value = 1.036083874717027
$ws.Range("K13").Value = 1.035660152375609
$ws.Range("L13").Value = 1.036898053806423
$ws.Range("M13").Value = 1.028931863614256
$ws.Range("N13").Value = 1.037555232659607
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029877032034797
$ws.Range("D14").Value = 1.032404838925001
$ws.Range("E14").Value = 1.033745563793204
$ws.Range("F14").Value = 1.02587403625041
$ws.Range("I14").Value = 1.03051581937658
$ws.Range("J14").Value = 1.036489499377522
$ws.Range("K14").Value = 1.035989057401952
$ws.Range("L14").Value = 1.037324775218897
$ws.Range("M14").Value = 1.029482851911781
$ws.Range("N14").Value = 1.037961433353646
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030165782055762
$ws.Range("D15").Value = 1.032627575879769
$ws.Range("E15").Value = 1.034028725690646
$ws.Range("F15").Value = 1.026234177797006
$ws.Range("I15").Value = 1.030591109347875
$ws.Range("J15").Value = 1.036739066343671
$ws.Range("K15").Value = 1.036191389356375
$ws.Range("L15").Value = 1.037587339723026
$ws.Range("M15").Value = 1.029821915706264
$ws.Range("N15").Value = 1.038211354733512
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03184095598537
$ws.Range("D16").Value = 1.033919312903821
$ws.Range("E16").Value = 1.035671577863044
$ws.Range("F16").Value = 1.028324022936046
$ws.Range("I16").Value = 1.031026748937297
$ws.Range("J16").Value = 1.038186347657936
$ws.Range("K16").Value = 1.037364257439719
$ws.Range("L16").Value = 1.039110261334059
$ws.Range("M16").Value = 1.031789117875249
$ws.Range("N16").Value = 1.039660691353249
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032887003580402
$ws.Range("D17").Value = 1.034725515796904
$ws.Range("E17").Value = 1.036697527388971
$ws.Range("F17").Value = 1.029629451270752
$ws.Range("I17").Value = 1.031297764511182
$ws.Range("J17").Value = 1.039089585044418
$ws.Range("K17").Value = 1.038095803741651
$ws.Range("L17").Value = 1.040060938611225
$ws.Range("M17").Value = 1.033017644452582
$ws.Range("N17").Value = 1.040565211440422
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03349546606197
$ws.Range("D18").Value = 1.035194316965717
$ws.Range("E18").Value = 1.037294329546461
$ws.Range("F18").Value = 1.030388952034586
$ws.Range("I18").Value = 1.031455041646086
$ws.Range("J18").Value = 1.039614796111753
$ws.Range("K18").Value = 1.038521024398469
$ws.Range("L18").Value = 1.040613818886964
$ws.Range("M18").Value = 1.033732298156617
$ws.Range("N18").Value = 1.041091168367715
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033702654044382
$ws.Range("D19").Value = 1.035353923408347
$ws.Range("E19").Value = 1.037497552167233
$ws.Range("F19").Value = 1.030647597954585
$ws.Range("I19").Value = 1.031508533955124
$ws.Range("J19").Value = 1.039793605231343
$ws.Range("K19").Value = 1.038665765054752
$ws.Range("L19").Value = 1.040802062295137
$ws.Range("M19").Value = 1.033975653632618
$ws.Range("N19").Value = 1.041270231416765
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.03277494696135
$ws.Range("D20").Value = 1.034639167667276
$ws.Range("E20").Value = 1.036587620609816
$ws.Range("F20").Value = 1.029489592049696
$ws.Range("I20").Value = 1.031268770227222
$ws.Range("J20").Value = 1.038992845634651
$ws.Range("K20").Value = 1.038017469113131
$ws.Range("L20").Value = 1.03995910954467
$ws.Range("M20").Value = 1.032886035145441
$ws.Range("N20").Value = 1.040468334649597
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029738928822292
$ws.Range("D21").Value = 1.032298300220937
$ws.Range("E21").Value = 1.0336101350059
$ws.Range("F21").Value = 1.025701796613744
$ws.Range("I21").Value = 1.03047978933222
$ws.Range("J21").Value = 1.036370126463665
$ws.Range("K21").Value = 1.035892269295825
$ws.Range("L21").Value = 1.03719919002346
$ws.Range("M21").Value = 1.029320686947342
$ws.Range("N21").Value = 1.037841890916558
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027814273499857
$ws.Range("D22").Value = 1.030812997344083
$ws.Range("E22").Value = 1.031722868055388
$ws.Range("F22").Value = 1.023301973574133
$ws.Range("I22").Value = 1.029976307418317
$ws.Range("J22").Value = 1.03470582838624
$ws.Range("K22").Value = 1.034542274831939
$ws.Range("L22").Value = 1.035448593164218
$ws.Range("M22").Value = 1.027060844442706
$ws.Range("N22").Value = 1.036175229344966
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028836099006942
$ws.Range("D23").Value = 1.031601689601525
$ws.Range("E23").Value = 1.032724816526432
$ws.Range("F23").Value = 1.024575941203995
$ws.Range("I23").Value = 1.030243925549349
$ws.Range("J23").Value = 1.035589582617387
$ws.Range("K23").Value = 1.03525926402156
$ws.Range("L23").Value = 1.036378100378372
$ws.Range("M23").Value = 1.028260592321616
$ws.Range("N23").Value = 1.037060238608489
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032825585708024
$ws.Range("D24").Value = 1.034678189126143
$ws.Range("E24").Value = 1.036637287743344
$ws.Range("F24").Value = 1.02955279437884
$ws.Range("I24").Value = 1.03128187397168
$ws.Range("J24").Value = 1.039036563048006
$ws.Range("K24").Value = 1.038052869722789
$ws.Range("L24").Value = 1.040005126759872
$ws.Range("M24").Value = 1.032945509667904
$ws.Range("N24").Value = 1.040512114146694
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037389601135783
$ws.Range("D25").Value = 1.038191917830773
$ws.Range("E25").Value = 1.041114373159042
$ws.Range("F25").Value = 1.035252723581862
$ws.Range("I25").Value = 1.032455001612292
$ws.Range("J25").Value = 1.042972852966841
$ws.Range("K25").Value = 1.041236948921693
$ws.Range("L25").Value = 1.044150294009111
$ws.Range("M25").Value = 1.038306973793555
$ws.Range("N25").Value = 1.044453994048713
